$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be auto-converted to numbers
$numericTextCells = @('D5', 'D6', 'D8', 'D9', 'D11', 'D12', 'D15', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D41', 'D42', 'D43', 'D44', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '27.366.18'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '1.712.37'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '224.78'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').Value = '0.5303'
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = '0.06684'
$ws.Range('E8').Value = '  +1.16%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '0.2669'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  -3.88%  '
$ws.Range('D11').Value = '0.07686'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '4.516'
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.716.68'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.948.16'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').Value = '0.5837'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '0.0₅8217'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').Value = '68.06'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '27.386.66'
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('D19').Value = '221.43'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').Value = '4.640'
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('D22').Value = '10.44'
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('D23').Value = '6.003'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('D24').Value = '1.005'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '144.70'
$ws.Range('E25').Value = '  -2.44%  '
$ws.Range('E26').Value = '  -2.62%  '
$ws.Range('E27').Value = '  -1.82%  '
$ws.Range('D28').Value = '7.266'
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('D29').Value = '16.23'
$ws.Range('E29').Value = '  -2.59%  '
$ws.Range('D30').Value = '0.05364'
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('D31').Value = '1.295'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').Value = '3.455'
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('D33').Value = '3.436'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('D34').Value = '1.642'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('D35').Value = '2.867'
$ws.Range('E35').Value = '  +1.04%  '
$ws.Range('D36').Value = '0.9527'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').Value = '2.398'
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('D38').Value = '0.5861'
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('D39').Value = '0.01639'
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('D40').Value = '1.093.55'
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range('D41').Value = '5.813'
$ws.Range('E41').Value = '  -1.92%  '
$ws.Range('D42').Value = '0.8448'
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('D43').Value = '1.005'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '101.02'
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('D45').Value = '1.855.13'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('D47').Value = '57.91'
$ws.Range('E47').Value = '  -1.99%  '
$ws.Range('D48').Value = '0.4533'
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('D49').Value = '1.005'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').Value = '8.088'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').Value = '0.05236'
$ws.Range('E51').Value = '  -0.32%  '

# Restore default style on the cells we temporarily reformatted
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).Style = "Normal"
}
